# Regenerate Handback status report timestamps.
# zh-cn sheet: handoff datetime 18:19:24 -> 18:20:46, handback datetime 18:19:48 -> 18:21:17
# de-de sheet: handoff datetime 18:19:29 -> 18:20:50, handback datetime 18:19:54 -> 18:21:24

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 18:20:46"
$wsZhCn.Range("E4").Value = "2016-03-22 18:20:46"
$wsZhCn.Range("H2").Value = "2016-03-22 18:21:17"
$wsZhCn.Range("H4").Value = "2016-03-22 18:21:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 18:20:50"
$wsDeDe.Range("E4").Value = "2016-03-22 18:20:50"
$wsDeDe.Range("H2").Value = "2016-03-22 18:21:24"
$wsDeDe.Range("H4").Value = "2016-03-22 18:21:24"
